# fixes to lca database based on discussion
# The "centralized" (YES/NO) column is replaced by a "scale" column
# (NONE/BUILDING/CITY/DISTRICT) on every sheet that has it.

$wb = $excel.ActiveWorkbook

# --- DHW sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("DHW")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Select()
$ws.Range("A16").Select()

# --- HEATING sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("HEATING")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Select()
$ws.Range("C8").Select()

# --- COOLING sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("COOLING")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "BUILDING"
$ws.Range("D5").Value = "DISTRICT"
$ws.Range("D6").Value = "DISTRICT"
$ws.Range("D7").Value = "DISTRICT"
$ws.Select()
$ws.Range("C11").Select()

# --- FUELS sheet ---------------------------------------------------------
# (no column value changes on this sheet, only selection / tab state)
$ws = $wb.Worksheets.Item("FUELS")
$ws.Select()
$ws.Range("C16").Select()

# --- ELECTRICITY sheet ----------------------------------------------------
# Activated/selected last so it becomes the workbook's active tab.
$ws = $wb.Worksheets.Item("ELECTRICITY")
$ws.Range("D1").Value = "scale"
$ws.Range("D2").Value = "NONE"
$ws.Range("D3").Value = "BUILDING"
$ws.Range("D4").Value = "CITY"
$ws.Select()
$ws.Range("D4").Select()
